$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1. Copy the formatting (style) from the
# existing H1 header cell so the new headers share the same cell style,
# then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the new I/J data columns for rows 2-12 (column I = "I0",
# column J = "IF").
$values = @{
    2  = @(8, 8)
    3  = @(7, 7)
    4  = @(6, 6)
    5  = @(8, 8)
    6  = @(9, 9)
    7  = @(8, 8)
    8  = @(4, 5)
    9  = @(3, 5)
    10 = @(7, 8)
    11 = @(3, 3)
    12 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
